$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item("Include ValueSets").Name = "Include ValueSet #0"
$wb.Worksheets.Item("Include from EntityCode").Name = "Include #1"

# --- Update Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row after "Contact" (row 10), before "Description" (row 11),
# copying row 10's formatting down so the new row keeps the same style.
$ws.Rows.Item(11).Insert()
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Fill new row 11 with Jurisdiction / (empty text value, not a truly blank cell)
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = "'"
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)

# Update Version value (row 3, column B)
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"
